# Apply updates to speech_noMFCCs_vowels_rfe.xlsx (sheet1 / ActiveSheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: Validation -> F1 train
$ws.Range("O1").Value = "F1 train"

# Row 2
$ws.Range("O2").Value = 0.8

# Row 3
$ws.Range("O3").Value = 1

# Row 4
$ws.Range("O4").Value = 1

# Row 5
$ws.Range("O5").Value = 1

# Row 6 (MLP, Technique 5) - parameters, TP, FN, and metrics updated
$ws.Range("C6").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("E6").Value = 9
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 0.85
$ws.Range("J6").Value = 0.8571428571428571
$ws.Range("K6").Value = 0.9
$ws.Range("L6").Value = 0.8181818181818182
$ws.Range("N6").Value = 0.9
$ws.Range("O6").Value = 0.7058823529411765

# Row 7
$ws.Range("O7").Value = 1

# Row 8
$ws.Range("O8").Value = 1

# Row 9
$ws.Range("O9").Value = 1

# Row 10
$ws.Range("O10").Value = 1

# Row 11 (MLP, Technique 10%) - parameters, TP, FN, and metrics updated
$ws.Range("C11").Value = "{'activation': 'relu', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("E11").Value = 9
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 0.65
$ws.Range("J11").Value = 0.72
$ws.Range("K11").Value = 0.9
$ws.Range("L11").Value = 0.6
$ws.Range("N11").Value = 0.9
$ws.Range("O11").Value = 0.7555555555555555

# Row 12
$ws.Range("O12").Value = 0.96

# Row 13
$ws.Range("O13").Value = 1

# Row 14
$ws.Range("O14").Value = 1

# Row 15
$ws.Range("O15").Value = 1

# Row 16
$ws.Range("O16").Value = 0.574468085106383
